$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 600
$ws.Range("I40").Value = 600
$ws.Range("K40").Value = 600
$ws.Range("M40").Value = -425
$ws.Range("H64").Value = 7551.4287
$ws.Range("I64").Value = 5981.8335
$ws.Range("K64").Value = 5981.8335
$ws.Range("M64").Value = -5733.8335
$ws.Range("H67").Value = 7551.4287
$ws.Range("I67").Value = 5981.8335
$ws.Range("K67").Value = 5981.8335
$ws.Range("M67").Value = -5123.8335
$ws.Range("H137").Value = 1593983.4
$ws.Range("I137").Value = 5143.4316
$ws.Range("K137").Value = 15430.2948
$ws.Range("M137").Value = -12880.2948
$ws.Range("H138").Value = 8911.296
$ws.Range("I138").Value = 15560.5
$ws.Range("J138").Value = 5000
$ws.Range("K138").Value = 46681.5
$ws.Range("L138").Value = 15000
$ws.Range("M138").Value = -41541.5
$ws.Range("N138").Value = -25280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 862016.9399999999
$ws.Range("J61").Value = 3329470
$ws.Range("L61").Value = 3329470
$ws.Range("N61").Value = -3329894
$ws.Range("H122").Value = 2073.1428
$ws.Range("I122").Value = 2073.1428
$ws.Range("K122").Value = 6219.428400000001
$ws.Range("M122").Value = -3769.428400000001
$ws.Range("H132").Value = 4063
$ws.Range("I132").Value = 2718.8
$ws.Range("J132").Value = 4735.1
$ws.Range("K132").Value = 8156.400000000001
$ws.Range("L132").Value = 14205.3
$ws.Range("M132").Value = -5626.400000000001
$ws.Range("N132").Value = -19265.3
$ws.Range("H136").Value = 862016.9399999999
$ws.Range("J136").Value = 3329470
$ws.Range("L136").Value = 9988410
$ws.Range("N136").Value = -9993510
$ws.Range("H138").Value = 72000
$ws.Range("J138").Value = 72000
$ws.Range("L138").Value = 72000
$ws.Range("N138").Value = -82280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5471.4287
$ws.Range("I86").Value = 3486.4285
$ws.Range("J86").Value = 7456.4287
$ws.Range("K86").Value = 3486.4285
$ws.Range("L86").Value = 7456.4287
$ws.Range("M86").Value = -2363.4285
$ws.Range("N86").Value = -9702.4287
$ws.Range("H89").Value = 5471.4287
$ws.Range("I89").Value = 3486.4285
$ws.Range("J89").Value = 7456.4287
$ws.Range("K89").Value = 17432.1425
$ws.Range("L89").Value = 37282.14350000001
$ws.Range("M89").Value = -11816.1425
$ws.Range("N89").Value = -48514.14350000001
$ws.Range("H94").Value = 2087.3667
$ws.Range("I94").Value = 1300.8125
$ws.Range("K94").Value = 1300.8125
$ws.Range("M94").Value = -849.8125
$ws.Range("H141").Value = 44370
$ws.Range("J141").Value = 49826.668
$ws.Range("L141").Value = 49826.668
$ws.Range("N141").Value = -60186.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1934.4
$ws.Range("I16").Value = 1934.4
$ws.Range("K16").Value = 1934.4
$ws.Range("M16").Value = -1647.4
$ws.Range("H31").Value = 2266.84
$ws.Range("I31").Value = 2268.1292
$ws.Range("J31").Value = 2264.7368
$ws.Range("K31").Value = 2268.1292
$ws.Range("L31").Value = 2264.7368
$ws.Range("M31").Value = -1973.1292
$ws.Range("N31").Value = -2854.7368
$ws.Range("H34").Value = 2266.84
$ws.Range("I34").Value = 2268.1292
$ws.Range("J34").Value = 2264.7368
$ws.Range("K34").Value = 2268.1292
$ws.Range("L34").Value = 2264.7368
$ws.Range("M34").Value = -2066.1292
$ws.Range("N34").Value = -2668.7368
$ws.Range("H58").Value = 1743.0769
$ws.Range("I58").Value = 1874.7
$ws.Range("J58").Value = 1604.5264
$ws.Range("K58").Value = 1874.7
$ws.Range("L58").Value = 1604.5264
$ws.Range("M58").Value = -1671.7
$ws.Range("N58").Value = -2010.5264
$ws.Range("H105").Value = 3435.7144
$ws.Range("I105").Value = 2950
$ws.Range("K105").Value = 2950
$ws.Range("M105").Value = -1203
$ws.Range("H113").Value = 1934.4
$ws.Range("I113").Value = 1934.4
$ws.Range("K113").Value = 1934.4
$ws.Range("M113").Value = 235.5999999999999
$ws.Range("H132").Value = 2975.0908
$ws.Range("I132").Value = 2165
$ws.Range("K132").Value = 6495
$ws.Range("M132").Value = -3965
$ws.Range("H134").Value = 2527.3
$ws.Range("I134").Value = 2344.238
$ws.Range("K134").Value = 7032.714
$ws.Range("M134").Value = -4497.714
$ws.Range("H136").Value = 1743.0769
$ws.Range("I136").Value = 1874.7
$ws.Range("J136").Value = 1604.5264
$ws.Range("K136").Value = 5624.1
$ws.Range("L136").Value = 4813.5792
$ws.Range("M136").Value = -3074.1
$ws.Range("N136").Value = -9913.5792

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 59.5
$ws.Range("J36").Value = 59.5
$ws.Range("L36").Value = 178.5
$ws.Range("N36").Value = -516.5
$ws.Range("H68").Value = 1991.75
$ws.Range("I68").Value = 2092.5
$ws.Range("J68").Value = 1891
$ws.Range("K68").Value = 6277.5
$ws.Range("L68").Value = 5673
$ws.Range("M68").Value = -5466.5
$ws.Range("N68").Value = -7295
$ws.Range("H71").Value = 1991.75
$ws.Range("I71").Value = 2092.5
$ws.Range("J71").Value = 1891
$ws.Range("K71").Value = 18832.5
$ws.Range("L71").Value = 17019
$ws.Range("M71").Value = -14776.5
$ws.Range("N71").Value = -25131
$ws.Range("H81").Value = 27784264
$ws.Range("I81").Value = 166667380
$ws.Range("J81").Value = 7642.2
$ws.Range("K81").Value = 500002140
$ws.Range("L81").Value = 22926.6
$ws.Range("M81").Value = -500001017
$ws.Range("N81").Value = -25172.6
$ws.Range("H84").Value = 27784264
$ws.Range("I84").Value = 166667380
$ws.Range("J84").Value = 7642.2
$ws.Range("K84").Value = 1500006420
$ws.Range("L84").Value = 68779.8
$ws.Range("M84").Value = -1500000804
$ws.Range("N84").Value = -80011.8
$ws.Range("H108").Value = 111114940
$ws.Range("J108").Value = 5500
$ws.Range("L108").Value = 16500
$ws.Range("N108").Value = -22260
$ws.Range("H117").Value = 333333340
$ws.Range("J117").Value = 333333340
$ws.Range("L117").Value = 1000000020
$ws.Range("N117").Value = -1000006904
$ws.Range("H131").Value = 2676935
$ws.Range("J131").Value = 3414.2856
$ws.Range("L131").Value = 10242.8568
$ws.Range("N131").Value = -20322.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 22728858
$ws.Range("I102").Value = 25001496
$ws.Range("K102").Value = 25001496
$ws.Range("M102").Value = -24999874
$ws.Range("H122").Value = 1845.2222
$ws.Range("I122").Value = 1066
$ws.Range("K122").Value = 3198
$ws.Range("M122").Value = -748
$ws.Range("H132").Value = 1167314.8
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9049.700000000001
$ws.Range("I7").Value = 4338.3335
$ws.Range("K7").Value = 4338.3335
$ws.Range("M7").Value = -4226.3335
$ws.Range("H16").Value = 1561.3
$ws.Range("J16").Value = 1496
$ws.Range("L16").Value = 1496
$ws.Range("N16").Value = -1836
$ws.Range("H25").Value = 10000
$ws.Range("J25").Value = 10000
$ws.Range("L25").Value = 10000
$ws.Range("N25").Value = -10460
$ws.Range("H40").Value = 3707227.8
$ws.Range("J40").Value = 3852.2
$ws.Range("L40").Value = 3852.2
$ws.Range("N40").Value = -4124.2
$ws.Range("H46").Value = 13625.909
$ws.Range("J46").Value = 5898.8
$ws.Range("L46").Value = 5898.8
$ws.Range("N46").Value = -6274.8
$ws.Range("H61").Value = 1271.48
$ws.Range("I61").Value = 1355.3684
$ws.Range("K61").Value = 1355.3684
$ws.Range("M61").Value = -1153.3684
$ws.Range("H104").Value = 26666
$ws.Range("J104").Value = 26666
$ws.Range("L104").Value = 26666
$ws.Range("N104").Value = -33654
$ws.Range("H113").Value = 1271.48
$ws.Range("I113").Value = 1355.3684
$ws.Range("K113").Value = 1355.3684
$ws.Range("M113").Value = 814.6315999999999
$ws.Range("H122").Value = 3250.5356
$ws.Range("I122").Value = 2771.739
$ws.Range("K122").Value = 8315.217000000001
$ws.Range("M122").Value = -5865.217000000001
$ws.Range("H126").Value = 9049.700000000001
$ws.Range("I126").Value = 4338.3335
$ws.Range("K126").Value = 13015.0005
$ws.Range("M126").Value = -10545.0005
$ws.Range("H132").Value = 4374.4
$ws.Range("I132").Value = 4674.8184
$ws.Range("J132").Value = 3548.25
$ws.Range("K132").Value = 14024.4552
$ws.Range("L132").Value = 10644.75
$ws.Range("M132").Value = -11494.4552
$ws.Range("N132").Value = -15704.75
$ws.Range("H136").Value = 28054.23
$ws.Range("J136").Value = 2826.3
$ws.Range("L136").Value = 8478.900000000001
$ws.Range("N136").Value = -13578.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").ClearContents()
$ws.Range("H49").Value = 99999
$ws.Range("I49").Value = 99999
$ws.Range("K49").Value = 99999
$ws.Range("M49").Value = -99769
$ws.Range("H122").Value = 1775.375
$ws.Range("I122").Value = 1671.8572
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 5015.571599999999
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -2565.571599999999
$ws.Range("N122").Value = -12400
$ws.Range("H132").Value = 28573436
$ws.Range("I132").Value = 47620436
$ws.Range("J132").Value = 2936.0715
$ws.Range("K132").Value = 142861308
$ws.Range("L132").Value = 8808.2145
$ws.Range("M132").Value = -142858778
$ws.Range("N132").Value = -13868.2145
